$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

# Update the daily conversion note (A1 on "Hoja1") with the new rates.
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.77 = 6410.05 pesos`n✅ 6410.05 pesos = 1.75 = 916.47 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Update the rate table on "tasas".
$ws2.Range("N10").Value = 566.453
$ws2.Range("O10").Value = 3630.99
$ws2.Range("N12").Value = 3665.01
$ws2.Range("O12").Value = 524
